$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily row (row 65) below the existing data (which ends at row 64)
$ws.Range("A65").Value = 46014
$ws.Range("B65").Value = 141
$ws.Range("C65").Value = 156
$ws.Range("D65").Value = 146

# Match the date formatting style used by the other cells in column A
$ws.Range("A65").NumberFormat = $ws.Range("A64").NumberFormat
